$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix sorting results: update Avg_Time_ms values for the affected rows.
$ws.Range("D2").Value = 0.6888
$ws.Range("D3").Value = 1.5352682

$wb.RefreshAll()
$excel.CalculateFullRebuild()
